$wb = $excel.ActiveWorkbook

# --- Sheet1: the irrigation-dates table rows below row 71 are wiped out ---
# (the user selected A72:B221 and cleared it entirely - contents AND formats -
# which is why the now-empty rows 124+ disappear completely from the sheet,
# while rows 72-123 survive as placeholders because column E still carries a
# styled, empty cell in those rows).
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A72:B221").Select()
$excel.Selection.Clear()

# The final selection/active-cell on Sheet1 after the edit.
$ws1.Range("A72:B221").Select()
$ws1.Range("A72").Activate()

# --- Sheet2: it is no longer the active/selected tab ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Select()

# Leave Sheet1 as the active sheet/tab.
$ws1.Activate()
